# Update countries & provincias Spain
#
# This script applies the data refresh captured in the commit "Update
# countries & provincias Spain". The source table (sheet "Pais") lists one
# country per row, sorted (roughly) by total cases. A few countries moved
# position in the ranking (Argentina, Uganda, Guyana and Belice each moved
# up past their neighbour), and a handful of rows received refreshed
# totals. The net effect, row by row, is captured explicitly below: for
# every affected row we set the country name (column A) and the 7 numeric
# columns (B:H = Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) to their final values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row => (Pais, CasosTotales, NuevosCasos, CasosActivos, Recuperados, CasosCriticos, MuertesHoy, Muertes)
$rows = @{
    4   = @("Estados Unidos", 1365357, 18048, 240865, 1043773, 16493, 682, 80719)
    15  = @("Canada", 68848, 1146, 31902, 32076, 502, 177, 4870)
    23  = @("Suiza", 30305, 54, 26600, 1872, 101, 3, 1833)
    52  = @("Noruega", 8105, 6, 32, 7854, 22, 0, 219)
    56  = @("Argentina", 6021, 245, 1757, 3959, 148, 5, 305)
    57  = @("Finlandia", 5962, 82, 4000, 1695, 45, 2, 267)
    157 = @("Uganda", 121, 5, 55, 66, 0, 0, 0)
    158 = @("Sudan del Sur", 120, 0, 2, 118, 0, 0, 0)
    159 = @("Bermudas", 118, 0, 64, 47, 4, 0, 7)
    162 = @("Guyana", 104, 10, 35, 59, 6, 0, 10)
    163 = @("Aruba", 101, 0, 89, 9, 4, 0, 3)
    192 = @("Belice", 18, 0, 16, 0, 0, 0, 2)
    193 = @("Nueva Caledonia", 18, 0, 18, 0, 0, 0, 0)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = $vals[6]
    $ws.Cells.Item($r, 8).Value = $vals[7]
}

$wb.Save()
